# Final upload edit: correct a handful of contact-detail values, re-blacken
# the pincode/phone number font color, and bump the header/data row heights.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (first contact) ---------------------------------------------
$ws.Range("B2").Value = "rohan"              # FirstName: tinu -> rohan
$ws.Range("J2").Value = "tintu@gmail.com"    # Email: tinu@gmail.com -> tintu@gmail.com
$ws.Range("L2").Value = "Reading ,Drawing"   # Hobbies: drop trailing comma

# --- Row 3 (second contact) ---------------------------------------------
$ws.Range("B3").Value = "mini"               # LastName: Maya -> mini
$ws.Range("H3").Value = "abcd"               # Street: dfbdf -> abcd
$ws.Range("L3").Value = "Reading ,Writing"   # Hobbies: drop trailing comma

# --- Formatting: make the Pincode/Phone numbers solid black -------------
$ws.Range("I2").Font.Color = 0
$ws.Range("K2").Font.Color = 0
$ws.Range("I3").Font.Color = 0
$ws.Range("K3").Font.Color = 0

# --- Row heights: header + both data rows grow slightly -----------------
$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 19.5
